$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Signature proportions")

# Insert a new column at the left; existing patient-code columns (and their
# data) shift from A:W to B:X.
$ws.Columns.Item(1).Insert()

# Header for the newly inserted column.
$ws.Range("A1").Value = "topic"

# Row-label ("topic") values for the 19 signature rows.
$topics = @("Age","SBS5","SBS8","SBS40","POLH","MMRD1","MMRD2","APOBEC1","APOBEC2","HRD","S-Dup","M-Dup","L-Dup","S-Del","L-Del","Clust-FBI","Clust-SV","Tr","FBI/Inv")

for ($i = 0; $i -lt $topics.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $topics[$i]
}
